# WAT API new test cases
# Adds three new rows (52-54) to the WoS_AuthorTransformation sheet describing
# new "Search Author API" test cases (WAT-805, WAT-806, WAT-807), matching the
# existing table layout/formatting, and moves the active selection to the new
# last row (A54) the way the author left it after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- seed row formatting ----------------------------------------------------
# Row 45 already carries the exact "blank STORE / blank DEPENDENCYTESTS" look
# (borders + top-vertical alignment + wrap) that rows 52 & 53 need, and row 54
# only differs in the STORE (K) cell, which mirrors row 10's populated-STORE
# style. Cloning formats from those existing rows keeps the new rows visually
# consistent with the rest of the sheet.
$ws.Range("A45:L45").Copy() | Out-Null
$ws.Range("A52:L52").PasteSpecial(-4122) | Out-Null
$ws.Range("A53:L53").PasteSpecial(-4122) | Out-Null
$ws.Range("A54:L54").PasteSpecial(-4122) | Out-Null

$ws.Range("K10").Copy() | Out-Null
$ws.Range("K54").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- row heights -------------------------------------------------------------
$ws.Rows("52:52").RowHeight = 30
$ws.Rows("54:54").RowHeight = 105

# --- row 52: WAT-805 ----------------------------------------------------------
$ws.Range("A52").Value = "WAT-805"
$ws.Range("B52").Value = 'Verify that Search Author API should return 400 if there is no "name" query param'
$ws.Range("C52").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D52").Value = "/author/search"
$ws.Range("E52").Value = "GET"
$ws.Range("G52").Value = "?&offset=0&limit=10`n"
$ws.Range("J52").Value = "status=400||error=Required query param 'name' is missing."

# --- row 53: WAT-806 ----------------------------------------------------------
$ws.Range("A53").Value = "WAT-806"
$ws.Range("B53").Value = "Verify that Search Author API should return 400 if there is empty name in query param"
$ws.Range("C53").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D53").Value = "/author/search"
$ws.Range("E53").Value = "GET"
$ws.Range("G53").Value = "?name=&offset=0&limit=10"
$ws.Range("J53").Value = 'status=400||error="A query param ''name'' is empty.'

# --- row 54: WAT-807 ----------------------------------------------------------
$ws.Range("A54").Value = "WAT-807"
$ws.Range("B54").Value = 'Verify that Search Author API should return results if there is atleast one non-blank "name" query param and ignore any "name" query parameter that is empty'
$ws.Range("C54").Value = "WOSAUTHORRECOMMEND"
$ws.Range("D54").Value = "/author/search"
$ws.Range("E54").Value = "GET"
$ws.Range("G54").Value = "?name=upadhyaya&offset=0&limit=10&name="
$ws.Range("J54").Value = "status=200||hits.primaryName=upadhyaya"
$ws.Range("K54").Value = "hits[0].authorClusterId||hits[0].primaryName||hits[0].alternativeName||hits[0].organization||hits[0].department||hits[0].location||hits[0].totalNumberOfPublications||hits[0].publicationYearRangeMin||hits[0].publicationYearRangeMax||hits[0].topPublications"

# --- leave the selection / scroll position where the author left it ----------
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A54").Select() | Out-Null
